{"js": "// The author re-typed the \"SOLOMON: Huzzah! ... The bag is filled only ...\"\n// line(s), which (per the XML diff) nets out to a single inserted word:\n// \"filled only the highest-quality...\" -> \"filled only with the highest-quality...\"\n// All surrounding paragraphs are unchanged. We locate the exact phrase and\n// insert the missing word \"with \" right after \"only \".\n\nconst searchText = \"filled only the highest-quality\";\nconst results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target phrase not found: \" + searchText);\n}\n\nconst target = results.items[0];\ntarget.insertText(\"filled only with the highest-quality\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The author re-typed the \"SOLOMON: Huzzah! ... The bag is filled only ...\"\n# line(s), which (per the XML diff) nets out to a single inserted word:\n# \"filled only the highest-quality...\" -> \"filled only with the highest-quality...\"\n# All surrounding paragraphs are unchanged. We locate the exact phrase and\n# insert the missing word \"with \" right after \"only \".\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"filled only the highest-quality\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"filled only with the highest-quality\"\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n"}
